$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 336.53845
$ws.Cells.Item(53, 9).Value = 410.85715
$ws.Cells.Item(53, 10).Value = 249.83333
$ws.Cells.Item(53, 11).Value = 410.85715
$ws.Cells.Item(53, 12).Value = 249.83333
$ws.Cells.Item(53, 13).Value = 226.14285
$ws.Cells.Item(53, 14).Value = -1523.83333

$ws.Cells.Item(107, 8).Value = 66.666664
$ws.Cells.Item(107, 9).Value = 66.666664
$ws.Cells.Item(107, 11).Value = 66.666664
$ws.Cells.Item(107, 13).Value = 1853.333336

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 5
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()

$ws.Cells.Item(11, 8).Value = 1599.5
$ws.Cells.Item(11, 9).Value = 200
$ws.Cells.Item(11, 10).Value = 2999
$ws.Cells.Item(11, 11).Value = 200
$ws.Cells.Item(11, 12).Value = 2999
$ws.Cells.Item(11, 13).Value = -56
$ws.Cells.Item(11, 14).Value = -3287

$ws.Cells.Item(45, 8).Value = 2140.3
$ws.Cells.Item(45, 9).Value = 1988
$ws.Cells.Item(45, 11).Value = 1988
$ws.Cells.Item(45, 13).Value = -1611

$ws.Cells.Item(102, 8).Value = 1755.5714
$ws.Cells.Item(102, 9).Value = 1481.5
$ws.Cells.Item(102, 10).Value = 3400
$ws.Cells.Item(102, 11).Value = 1481.5
$ws.Cells.Item(102, 12).Value = 3400
$ws.Cells.Item(102, 13).Value = 140.5
$ws.Cells.Item(102, 14).Value = -6644

$ws.Cells.Item(106, 8).Value = 20000
$ws.Cells.Item(106, 10).Value = 20000
$ws.Cells.Item(106, 12).Value = 20000
$ws.Cells.Item(106, 14).Value = -22524

$ws.Cells.Item(110, 8).Value = 819.2
$ws.Cells.Item(110, 9).Value = 699.1667
$ws.Cells.Item(110, 10).Value = 999.25
$ws.Cells.Item(110, 11).Value = 699.1667
$ws.Cells.Item(110, 12).Value = 999.25
$ws.Cells.Item(110, 13).Value = 1345.8333
$ws.Cells.Item(110, 14).Value = -5089.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3369.5
$ws.Cells.Item(20, 9).Value = 3428.1428
$ws.Cells.Item(20, 11).Value = 3428.1428
$ws.Cells.Item(20, 13).Value = -3181.1428

$ws.Cells.Item(35, 8).Value = 14982.5
$ws.Cells.Item(35, 10).Value = 14982.5
$ws.Cells.Item(35, 12).Value = 14982.5
$ws.Cells.Item(35, 14).Value = -15602.5

$ws.Cells.Item(80, 8).Value = 540.5
$ws.Cells.Item(80, 9).Value = 384.75
$ws.Cells.Item(80, 10).Value = 644.3333
$ws.Cells.Item(80, 11).Value = 384.75
$ws.Cells.Item(80, 12).Value = 644.3333
$ws.Cells.Item(80, 13).Value = 613.25
$ws.Cells.Item(80, 14).Value = -2640.3333

$ws.Cells.Item(83, 8).Value = 540.5
$ws.Cells.Item(83, 9).Value = 384.75
$ws.Cells.Item(83, 10).Value = 644.3333
$ws.Cells.Item(83, 11).Value = 1923.75
$ws.Cells.Item(83, 12).Value = 3221.6665
$ws.Cells.Item(83, 13).Value = 3068.25
$ws.Cells.Item(83, 14).Value = -13205.6665

$ws.Cells.Item(105, 8).Value = 4242.875
$ws.Cells.Item(105, 9).Value = 4424.3335
$ws.Cells.Item(105, 10).Value = 3698.5
$ws.Cells.Item(105, 11).Value = 4424.3335
$ws.Cells.Item(105, 12).Value = 3698.5
$ws.Cells.Item(105, 13).Value = -2677.3335
$ws.Cells.Item(105, 14).Value = -7192.5

$ws.Cells.Item(107, 8).Value = 1999.3334
$ws.Cells.Item(107, 9).Value = 1999
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 1999
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = -79
$ws.Cells.Item(107, 14).Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(131, 8).Value = 49994.5
$ws.Cells.Item(131, 10).Value = 49994.5
$ws.Cells.Item(131, 12).Value = 49994.5
$ws.Cells.Item(131, 14).Value = -60074.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 522
$ws.Cells.Item(97, 10).Value = 582.2
$ws.Cells.Item(97, 12).Value = 1746.6
$ws.Cells.Item(97, 14).Value = -2738.6

$ws.Cells.Item(117, 8).Value = 548.8333
$ws.Cells.Item(117, 9).Value = 429
$ws.Cells.Item(117, 11).Value = 1287
$ws.Cells.Item(117, 13).Value = 2155

$ws.Cells.Item(140, 8).Value = 1500
$ws.Cells.Item(140, 9).Value = 1500
$ws.Cells.Item(140, 11).Value = 4500
$ws.Cells.Item(140, 13).Value = 680

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9078.799999999999
$ws.Cells.Item(70, 9).Value = 7997.5
$ws.Cells.Item(70, 11).Value = 7997.5
$ws.Cells.Item(70, 13).Value = -7727.5

$ws.Cells.Item(73, 8).Value = 9078.799999999999
$ws.Cells.Item(73, 9).Value = 7997.5
$ws.Cells.Item(73, 11).Value = 7997.5
$ws.Cells.Item(73, 13).Value = -7061.5

$ws.Cells.Item(80, 8).Value = 2748.3333
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 2748.3333
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 2748.3333
$ws.Cells.Item(80, 14).Value = -4744.3333
$ws.Cells.Item(80, 13).ClearContents()

$ws.Cells.Item(83, 8).Value = 2748.3333
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 2748.3333
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 13741.6665
$ws.Cells.Item(83, 14).Value = -23725.6665
$ws.Cells.Item(83, 13).ClearContents()

$ws.Cells.Item(97, 8).Value = 888.3333
$ws.Cells.Item(97, 9).Value = 883.1667
$ws.Cells.Item(97, 10).Value = 898.6667
$ws.Cells.Item(97, 11).Value = 883.1667
$ws.Cells.Item(97, 12).Value = 898.6667
$ws.Cells.Item(97, 13).Value = -387.1667
$ws.Cells.Item(97, 14).Value = -1890.6667

$ws.Cells.Item(107, 8).Value = 63
$ws.Cells.Item(107, 9).Value = 45
$ws.Cells.Item(107, 11).Value = 45
$ws.Cells.Item(107, 13).Value = 1875

$ws.Cells.Item(126, 8).Value = 13195.429
$ws.Cells.Item(126, 9).Value = 12303.818
$ws.Cells.Item(126, 10).Value = 16464.666
$ws.Cells.Item(126, 11).Value = 36911.454
$ws.Cells.Item(126, 12).Value = 49393.99800000001
$ws.Cells.Item(126, 13).Value = -34441.454
$ws.Cells.Item(126, 14).Value = -54333.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(10, 8).Value = 2500
$ws.Cells.Item(10, 9).Value = 2500
$ws.Cells.Item(10, 11).Value = 2500
$ws.Cells.Item(10, 13).Value = -2360

$ws.Cells.Item(55, 9).Value = 277.41666
$ws.Cells.Item(55, 10).Value = 277.07693
$ws.Cells.Item(55, 11).Value = 277.41666
$ws.Cells.Item(55, 12).Value = 277.07693
$ws.Cells.Item(55, 13).Value = -104.41666
$ws.Cells.Item(55, 14).Value = -623.0769299999999

$ws.Cells.Item(122, 8).Value = 3504
$ws.Cells.Item(122, 9).Value = 3504
$ws.Cells.Item(122, 11).Value = 10512
$ws.Cells.Item(122, 13).Value = -8062

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 13).ClearContents()

$ws.Cells.Item(11, 8).Value = 41667004
$ws.Cells.Item(11, 10).Value = 505
$ws.Cells.Item(11, 12).Value = 505
$ws.Cells.Item(11, 14).Value = -789

$ws.Cells.Item(107, 8).Value = 469.84616
$ws.Cells.Item(107, 9).Value = 326.83334
$ws.Cells.Item(107, 11).Value = 980.5000200000001
$ws.Cells.Item(107, 13).Value = 939.4999799999999
